# Afegir contingut taula i r
# Replace the remaining "TARDA" placeholder cells (rows 7-8, the Fly EM
# rows) with the real "SI" value, and flip the still-pending "NO" cells
# in rows 12-14 (Fly Opt - er/sw/ba) to "SI" to finish filling in the
# completion table. Also moves the sheet's active selection to where the
# author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 7 & 8 ("Fly EM" / "Fly EM - er"): "TARDA" -> "SI" across columns
# C, D, E, F, G, I (H, J, K, L, M stay "NO").
foreach ($r in 7, 8) {
    foreach ($col in "C", "D", "E", "F", "G", "I") {
        $ws.Range("$col$r").Value = "SI"
    }
}

# Row 12 ("Fly Opt - er"): "NO" -> "SI" for C, F, G, I.
$ws.Range("C12").Value = "SI"
$ws.Range("F12").Value = "SI"
$ws.Range("G12").Value = "SI"
$ws.Range("I12").Value = "SI"

# Row 13 ("Fly Opt - sw"): "NO" -> "SI" for F, I.
$ws.Range("F13").Value = "SI"
$ws.Range("I13").Value = "SI"

# Row 14 ("Fly Opt - ba"): "NO" -> "SI" for F, I.
$ws.Range("F14").Value = "SI"
$ws.Range("I14").Value = "SI"

# Scroll the view down a bit and move the active selection, matching
# where the author's cursor ended up after finishing the edits.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("I14").Select()
